# Add a "Status" column (C) to the build-tracking sheet, marking the
# items that are done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122

# Header cell C2 "Status" - copy formatting from B2, then set text.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C2").Value = "Status"

# Rows 3 & 4 (DW Exist Statement / DW Build Statement) -> Done
$ws.Range("B3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C3").Value = "Done"

$ws.Range("B4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C4").Value = "Done"

# Rows 6-10 (FactProudctionBatch, DimDate, DimShipment, DimCompliance, DimOrder) -> Done
$ws.Range("B6").Copy() | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C6").Value = "Done"

$ws.Range("B7").Copy() | Out-Null
$ws.Range("C7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C7").Value = "Done"

$ws.Range("B8").Copy() | Out-Null
$ws.Range("C8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C8").Value = "Done"

$ws.Range("B9").Copy() | Out-Null
$ws.Range("C9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C9").Value = "Done"

$ws.Range("B10").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").Value = "Done"

# Row 11 (DimItem) -> not yet done, only formatting carried over, no value.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("C11").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# Update the active cell selection to match the author's final state.
$ws.Range("C11").Select() | Out-Null
